# Update Wild Card round row (row 3, label "R") on both the OFF and DEF sheets
# with the simulated game's target depth stats.

$wb = $excel.ActiveWorkbook

$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 222
$wsOff.Range("C3").Value = 130
$wsOff.Range("D3").Value = 45
$wsOff.Range("E3").Value = 22

$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 185
$wsDef.Range("C3").Value = 130
$wsDef.Range("D3").Value = 64
$wsDef.Range("E3").Value = 30
